$wb = $excel.ActiveWorkbook

# --- Department sheet (sheet4): update test data rows 2-5, add new I11 cell ---
$wsDept = $wb.Worksheets.Item("Department")
$wsDept.Range("A2").Value = "testExcel542"
$wsDept.Range("A3").Value = "testExcel632"
$wsDept.Range("A4").Value = "testExcel712"
$wsDept.Range("A5").Value = "testExcel882"
$wsDept.Cells.Item(11, 9).NumberFormat = "mm-dd-yy"
$wsDept.Columns.Item(9).ColumnWidth = 9.666666667

# --- Affiliate sheet (sheet1): update test data rows 2-4 ---
$wsAff = $wb.Worksheets.Item("Affiliate")
$wsAff.Range("A2").Value = "testAffiliateexcel1211"
$wsAff.Range("A3").Value = "testAffiliateexcel2111"
$wsAff.Range("A4").Value = "testAffiliateexcel3311"
$wsAff.Range("B2").Value = "'1234564783901"
$wsAff.Range("B3").Value = "'12345657893901"
$wsAff.Range("B4").Value = "'12345667893001"

# --- Selections / active tab bookkeeping ---
# Department keeps a recorded selection at A5 (no longer the active tab)
$wsDept.Activate()
$wsDept.Range("A5").Select()

# Affiliate becomes the active tab with selection at B4
$wsAff.Activate()
$wsAff.Range("B4").Select()
